# The sheet originally had two extra columns - "linkage" (F) and
# "zygosity" (G) - inserted between "allele" (E) and "allele_count"
# (formerly H). Remove them so the remaining columns shift left,
# matching the upstream schema (A1:M2 instead of A1:O2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F:G").Delete()
